$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $val)
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '30.203.10'
Set-TextValue 'E2' '  -1.17%  '
Set-TextValue 'D3' '1.858.90'
Set-TextValue 'E3' '  -3.21%  '
Set-TextValue 'D4' '1.001'
Set-TextValue 'E4' '  +0.09%  '
Set-TextValue 'D5' '234.07'
Set-TextValue 'E5' '  -2.73%  '
Set-TextValue 'D6' '1.000'
Set-TextValue 'E6' '  +0.06%  '
Set-TextValue 'D7' '0.4682'
Set-TextValue 'E7' '  -2.24%  '
Set-TextValue 'E8' '  -1.32%  '
Set-TextValue 'E9' '  -2.77%  '
Set-TextValue 'D10' '20.21'
Set-TextValue 'E10' '  +3.55%  '
Set-TextValue 'D11' '0.07818'
Set-TextValue 'E11' '  +0.42%  '
Set-TextValue 'D12' '96.91'
Set-TextValue 'E12' '  -7.10%  '
Set-TextValue 'D13' '1.866.00'
Set-TextValue 'E13' '  -2.75%  '
Set-TextValue 'D14' '5.101'
Set-TextValue 'E14' '  -2.60%  '
Set-TextValue 'D15' '0.6646'
Set-TextValue 'E15' '  -1.64%  '
Set-TextValue 'D16' '283.31'
Set-TextValue 'E16' '  -3.44%  '
Set-TextValue 'D17' '30.255.22'
Set-TextValue 'E17' '  -0.98%  '
Set-TextValue 'D18' '1.000'
Set-TextValue 'E18' '  +0.11%  '
Set-TextValue 'D19' '5.436'
Set-TextValue 'E19' '  -0.46%  '
Set-TextValue 'D20' '12.61'
Set-TextValue 'E20' '  -1.35%  '
Set-TextValue 'D21' '2.108.61'
Set-TextValue 'E21' '  -2.56%  '
Set-TextValue 'D22' '0.000007245'
Set-TextValue 'E22' '  -3.64%  '
Set-TextValue 'D23' '0.9998'
Set-TextValue 'E23' '  -0.01%  '
Set-TextValue 'D24' '6.135'
Set-TextValue 'E24' '  -3.49%  '
Set-TextValue 'D25' '168.04'
Set-TextValue 'E25' '  +0.19%  '
Set-TextValue 'D26' '9.314'
Set-TextValue 'E26' '  -1.44%  '
Set-TextValue 'D27' '18.97'
Set-TextValue 'E27' '  -3.25%  '
Set-TextValue 'E28' '  -9.32%  '
Set-TextValue 'D29' '1.343'
Set-TextValue 'E29' '  -3.37%  '
Set-TextValue 'D30' '0.09587'
Set-TextValue 'E30' '  -3.87%  '
Set-TextValue 'D31' '4.393'
Set-TextValue 'E31' '  -4.43%  '
Set-TextValue 'D32' '1.470'
Set-TextValue 'E32' '  -3.25%  '
Set-TextValue 'D33' '4.099'
Set-TextValue 'E33' '  -4.47%  '
Set-TextValue 'D34' '0.04659'
Set-TextValue 'E34' '  -2.11%  '
Set-TextValue 'B35' 'ImmutableX'
Set-TextValue 'C35' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D35' '0.6994'
Set-TextValue 'E35' '  -4.31%  '
Set-TextValue 'B36' 'ARBITRUM'
Set-TextValue 'C36' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 'D36' '1.095'
Set-TextValue 'E36' '  -2.10%  '
Set-TextValue 'D37' '0.9995'
Set-TextValue 'E37' '  +0.14%  '
Set-TextValue 'D38' '2.709'
Set-TextValue 'E38' '  -0.09%  '
Set-TextValue 'D39' '0.01853'
Set-TextValue 'E39' '  -3.48%  '
Set-TextValue 'D40' '6.432'
Set-TextValue 'E40' '  -0.11%  '
Set-TextValue 'D41' '2.510'
Set-TextValue 'E41' '  -4.43%  '
Set-TextValue 'D42' '72.19'
Set-TextValue 'E42' '  -3.94%  '
Set-TextValue 'D43' '0.8542'
Set-TextValue 'E43' '  -1.01%  '
Set-TextValue 'D44' '1.931'
Set-TextValue 'E44' '  -2.09%  '
Set-TextValue 'D45' '104.15'
Set-TextValue 'E45' '  -1.98%  '
Set-TextValue 'D46' '0.4162'
Set-TextValue 'E46' '  -2.97%  '
Set-TextValue 'D47' '0.9997'
Set-TextValue 'E47' '  +0.00%  '
Set-TextValue 'D48' '1.007.68'
Set-TextValue 'E48' '  +2.69%  '
Set-TextValue 'D49' '7.190'
Set-TextValue 'E49' '  -3.64%  '
Set-TextValue 'D50' '9.072'
Set-TextValue 'E50' '  +2.25%  '
Set-TextValue 'D51' '33.80'
Set-TextValue 'E51' '  -2.98%  '
